$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on every touched cell before assigning, so Excel
# does not auto-coerce numeric-looking strings (e.g. "4.20", "0.197",
# "43.754.35") into numbers and silently drop formatting/precision.


# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.754.35'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.21%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.243.18'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.14%  '

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.19%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '318.04'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.92%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '100.94'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.28%  '

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -1.38%  '

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.11%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.554'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.87%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.85'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.06%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0830'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.28%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.64'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.33%  '

# Row 13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.88%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.588.33'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.37%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.854'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.80%  '

# Row 16
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.256.09'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.85%  '

# Row 17
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = 'Chainlink'
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.15'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.28%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.679.04'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.21%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.32'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -8.05%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0982'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.45%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.51'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.74%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '65.53'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.14%  '

# Row 23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -2.16%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '235.06'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -1.41%  '

# Row 25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -3.20%  '

# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.06%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.08'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.17%  '

# Row 28
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.19'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -2.65%  '

# Row 29
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = 'InjectiveProtocol'
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '37.03'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +2.72%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.22'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -2.82%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '158.39'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +3.48%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.09'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.85%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0847'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -4.25%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.69'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.81%  '

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +11.23%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.05'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -5.08%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.93'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -1.31%  '

# Row 38
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -3.04%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.73'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +1.25%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.20'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -5.68%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '15.79'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +16.60%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0315'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -3.38%  '

# Row 43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.27%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.785.61'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +1.26%  '

# Row 45
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = 'ordi'
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '75.27'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.81%  '

# Row 46
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = 'Algorand'
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.197'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -4.10%  '

# Row 47
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = 'BitcoinSV'
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '82.28'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -5.23%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.18'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -2.77%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '58.21'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.74%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '103.09'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.49%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.66'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +3.03%  '
